# Auto-generated script to update Zeromus_Profits market data values
# across sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 1154.9286
$ws.Range("I33").Value = 828.4231
$ws.Range("J33").Value = 5399.5
$ws.Range("K33").Value = 828.4231
$ws.Range("L33").Value = 5399.5
$ws.Range("M33").Value = -599.4231
$ws.Range("N33").Value = -5857.5
$ws.Range("H63").Value = 19999.334
$ws.Range("J63").Value = 19999.334
$ws.Range("L63").Value = 19999.334
$ws.Range("N63").Value = -21247.334
$ws.Range("H66").Value = 19999.334
$ws.Range("J66").Value = 19999.334
$ws.Range("L66").Value = 59998.00199999999
$ws.Range("N66").Value = -66238.00199999999
$ws.Range("H100").Value = 6109.737
$ws.Range("I100").Value = 1480.8334
$ws.Range("J100").Value = 8246.154
$ws.Range("K100").Value = 1480.8334
$ws.Range("L100").Value = 8246.154
$ws.Range("M100").Value = -939.8334
$ws.Range("N100").Value = -9328.154
$ws.Range("H125").Value = 1963.6666
$ws.Range("I125").Value = 856.4
$ws.Range("K125").Value = 7707.599999999999
$ws.Range("M125").Value = -5247.599999999999
$ws.Range("H129").Value = 1764.341
$ws.Range("I129").Value = 348.75
$ws.Range("J129").Value = 2078.9167
$ws.Range("K129").Value = 1046.25
$ws.Range("L129").Value = 6236.750100000001
$ws.Range("M129").Value = 3953.75
$ws.Range("N129").Value = -16236.7501
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10059.044
$ws.Range("I32").Value = 2867.4707
$ws.Range("J32").Value = 30435.166
$ws.Range("K32").Value = 2867.4707
$ws.Range("L32").Value = 30435.166
$ws.Range("M32").Value = -2580.4707
$ws.Range("N32").Value = -31009.166
$ws.Range("H61").Value = 1945.7826
$ws.Range("I61").Value = 1844.8948
$ws.Range("J61").Value = 2425
$ws.Range("K61").Value = 1844.8948
$ws.Range("L61").Value = 2425
$ws.Range("M61").Value = -1632.8948
$ws.Range("N61").Value = -2849
$ws.Range("H74").Value = 26474848
$ws.Range("I74").Value = 37504070
$ws.Range("J74").Value = 4705.6
$ws.Range("K74").Value = 37504070
$ws.Range("L74").Value = 4705.6
$ws.Range("M74").Value = -37503196
$ws.Range("N74").Value = -6453.6
$ws.Range("H77").Value = 26474848
$ws.Range("I77").Value = 37504070
$ws.Range("J77").Value = 4705.6
$ws.Range("K77").Value = 187520350
$ws.Range("L77").Value = 23528
$ws.Range("M77").Value = -187515982
$ws.Range("N77").Value = -32264
$ws.Range("H132").Value = 1927.4615
$ws.Range("I132").Value = 1644.8182
$ws.Range("J132").Value = 3482
$ws.Range("K132").Value = 4934.4546
$ws.Range("L132").Value = 10446
$ws.Range("M132").Value = -2404.4546
$ws.Range("N132").Value = -15506
$ws.Range("H136").Value = 1945.7826
$ws.Range("I136").Value = 1844.8948
$ws.Range("J136").Value = 2425
$ws.Range("K136").Value = 5534.6844
$ws.Range("L136").Value = 7275
$ws.Range("M136").Value = -2984.6844
$ws.Range("N136").Value = -12375
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1545.8368
$ws.Range("I20").Value = 1528.0571
$ws.Range("J20").Value = 1590.2858
$ws.Range("K20").Value = 1528.0571
$ws.Range("L20").Value = 1590.2858
$ws.Range("M20").Value = -1281.0571
$ws.Range("N20").Value = -2084.2858
$ws.Range("H94").Value = 7512.2856
$ws.Range("I94").Value = 372.6
$ws.Range("J94").Value = 25361.5
$ws.Range("K94").Value = 372.6
$ws.Range("L94").Value = 25361.5
$ws.Range("M94").Value = 78.39999999999998
$ws.Range("N94").Value = -26263.5
$ws.Range("H134").Value = 2768.6843
$ws.Range("I134").Value = 1940.3334
$ws.Range("K134").Value = 5821.0002
$ws.Range("M134").Value = -3286.0002
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 83334610
$ws.Range("I94").Value = 333334000
$ws.Range("J94").Value = 1479.3334
$ws.Range("K94").Value = 333334000
$ws.Range("L94").Value = 1479.3334
$ws.Range("M94").Value = -333333549
$ws.Range("N94").Value = -2381.3334
$ws.Range("H105").Value = 1729
$ws.Range("I105").Value = 1720
$ws.Range("J105").Value = 1765
$ws.Range("K105").Value = 1720
$ws.Range("L105").Value = 1765
$ws.Range("M105").Value = 27
$ws.Range("N105").Value = -5259
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 55.333332
$ws.Range("I12").Value = 11.666667
$ws.Range("J12").Value = 61.57143
$ws.Range("K12").Value = 35.000001
$ws.Range("L12").Value = 184.71429
$ws.Range("M12").Value = 137.999999
$ws.Range("N12").Value = -530.71429
$ws.Range("H107").Value = 523.4666999999999
$ws.Range("I107").Value = 284.66666
$ws.Range("J107").Value = 881.6667
$ws.Range("K107").Value = 853.9999799999999
$ws.Range("L107").Value = 2645.0001
$ws.Range("M107").Value = 1066.00002
$ws.Range("N107").Value = -6485.0001
$ws.Range("H131").Value = 1007.8372
$ws.Range("I131").Value = 742.1
$ws.Range("K131").Value = 2226.3
$ws.Range("M131").Value = 2813.7
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 26671746
$ws.Range("I70").Value = 47063930
$ws.Range("J70").Value = 5044.615
$ws.Range("K70").Value = 47063930
$ws.Range("L70").Value = 5044.615
$ws.Range("M70").Value = -47063660
$ws.Range("N70").Value = -5584.615
$ws.Range("H73").Value = 26671746
$ws.Range("I73").Value = 47063930
$ws.Range("J73").Value = 5044.615
$ws.Range("K73").Value = 47063930
$ws.Range("L73").Value = 5044.615
$ws.Range("M73").Value = -47062994
$ws.Range("N73").Value = -6916.615
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()
$ws.Range("H102").Value = 1862.3158
$ws.Range("I102").Value = 1907.6
$ws.Range("J102").Value = 1692.5
$ws.Range("K102").Value = 1907.6
$ws.Range("L102").Value = 1692.5
$ws.Range("M102").Value = -285.5999999999999
$ws.Range("N102").Value = -4936.5
$ws.Range("H107").Value = 33333974
$ws.Range("I107").Value = 83333630
$ws.Range("K107").Value = 83333630
$ws.Range("M107").Value = -83331710
$ws.Range("H132").Value = 3259.2
$ws.Range("I132").Value = 1932.3334
$ws.Range("K132").Value = 5797.0002
$ws.Range("M132").Value = -3267.0002
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1705.3334
$ws.Range("I7").Value = 1377.9166
$ws.Range("J7").Value = 2360.1667
$ws.Range("K7").Value = 1377.9166
$ws.Range("L7").Value = 2360.1667
$ws.Range("M7").Value = -1265.9166
$ws.Range("N7").Value = -2584.1667
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("M34").ClearContents()
$ws.Range("N34").ClearContents()
$ws.Range("H126").Value = 1705.3334
$ws.Range("I126").Value = 1377.9166
$ws.Range("J126").Value = 2360.1667
$ws.Range("K126").Value = 4133.7498
$ws.Range("L126").Value = 7080.500100000001
$ws.Range("M126").Value = -1663.7498
$ws.Range("N126").Value = -12020.5001
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 13272.5
$ws.Range("I51").Value = 9000
$ws.Range("K51").Value = 9000
$ws.Range("M51").Value = -8490
$ws.Range("H113").Value = 15625829
$ws.Range("I113").Value = 29412658
$ws.Range("J113").Value = 756.3333
$ws.Range("K113").Value = 88237974
$ws.Range("L113").Value = 2268.9999
$ws.Range("M113").Value = -88235804
$ws.Range("N113").Value = -6608.9999
$ws.Range("H123").Value = 30988.166
$ws.Range("J123").Value = 30988.166
$ws.Range("L123").Value = 30988.166
$ws.Range("N123").Value = -40788.166
$ws.Range("H126").Value = 3617.4666
$ws.Range("I126").Value = 3966.3076
$ws.Range("J126").Value = 1350
$ws.Range("K126").Value = 11898.9228
$ws.Range("L126").Value = 4050
$ws.Range("M126").Value = -9428.9228
$ws.Range("N126").Value = -8990
